$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.730.89'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.849.35'
$ws.Range("E3").Value = '  +0.12%  '
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.91'
$ws.Range("E5").Value = '  -0.16%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.002'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4323'
$ws.Range("E7").Value = '  +1.35%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3657'
$ws.Range("E8").Value = '  -0.51%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.14'
$ws.Range("E9").Value = '  +0.85%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07339'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.8792'
$ws.Range("E11").Value = '  -2.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '20.73'
$ws.Range("E12").Value = '  +0.13%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.856.67'
$ws.Range("E13").Value = '  +1.97%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.342'
$ws.Range("E14").Value = '  -0.24%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.526'
$ws.Range("E15").Value = '  -0.79%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.06947'
$ws.Range("E16").Value = '  +1.64%  '
$ws.Range("E17").Value = '  +0.00%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '80.29'
$ws.Range("E18").Value = '  +3.37%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000009049'
$ws.Range("E19").Value = '  +2.45%  '
$ws.Range("E20").Value = '  +0.24%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '15.39'
$ws.Range("E21").Value = '  -0.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '27.962.96'
$ws.Range("E22").Value = '  +1.13%  '
$ws.Range("E23").Value = '  +0.33%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '10.38'
$ws.Range("E24").Value = '  -2.68%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.089.99'
$ws.Range("E25").Value = '  +1.38%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.989'
$ws.Range("E26").Value = '  -2.86%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '155.73'
$ws.Range("E27").Value = '  +1.11%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.63'
$ws.Range("E28").Value = '  +2.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '120.69'
$ws.Range("E29").Value = '  +8.70%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.273'
$ws.Range("E30").Value = '  +0.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.863'
$ws.Range("E31").Value = '  +1.89%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.08927'
$ws.Range("E32").Value = '  +0.36%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7588'
$ws.Range("E33").Value = '  -1.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.552'
$ws.Range("E34").Value = '  -0.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.951'
$ws.Range("E35").Value = '  +1.16%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.122'
$ws.Range("E36").Value = '  +3.55%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.107'
$ws.Range("E37").Value = '  +1.11%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.05419'
$ws.Range("E38").Value = '  +0.73%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01937'
$ws.Range("E39").Value = '  +0.68%  '
$ws.Range("E40").Value = '  -3.64%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.5098'
$ws.Range("E41").Value = '  +0.66%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1660'
$ws.Range("E42").Value = '  +1.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.669'
$ws.Range("E43").Value = '  -1.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.344'
$ws.Range("E44").Value = '  +0.95%  '
$ws.Range("E45").Value = '  +1.41%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.06541'
$ws.Range("E46").Value = '  -1.51%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4665'
$ws.Range("E47").Value = '  -1.07%  '
$ws.Range("E48").Value = '  -0.73%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.001'
$ws.Range("E49").Value = '  -0.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.622'
$ws.Range("E50").Value = '  -1.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '64.73'
$ws.Range("E51").Value = '  +0.69%  '
